# Weekly fruit/hortaliza price update:
# Insert a new week's worth of data (2 rows: "Primera" and "Segunda" quality
# grades) for the "Acelga" series at Mercado Mayorista Lo Valledor de
# Santiago, right after the existing row for date serial 44265 (row 826),
# pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 827-828 (everything from the old row 827 onward
# shifts down by two, landing on rows 829-860; dimension grows to R860).
$ws.Range("A827:A828").EntireRow.Insert()

# New row 827: Acelga / Primera
$ws.Range("A827").Value = 6
$ws.Range("B827").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C827").Value = "Metropolitana"
$ws.Range("D827").Value = 44747
$ws.Range("E827").Value = 13
$ws.Range("F827").Value = 100112009
$ws.Range("G827").Value = "Acelga"
$ws.Range("H827").Value = "Sin especificar"
$ws.Range("I827").Value = "Primera"
$ws.Range("J827").Value = 90
$ws.Range("K827").Value = 16000
$ws.Range("L827").Value = 16000
$ws.Range("M827").Value = 16000
$ws.Range("N827").Value = "`$/docena de atados"
$ws.Range("O827").Value = "Región Metropolitana"
$ws.Range("P827").Value = 5333
$ws.Range("Q827").Value = 3
$ws.Range("R827").Value = "Hortaliza"

# New row 828: Acelga / Segunda
$ws.Range("A828").Value = 6
$ws.Range("B828").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C828").Value = "Metropolitana"
$ws.Range("D828").Value = 44747
$ws.Range("E828").Value = 13
$ws.Range("F828").Value = 100112009
$ws.Range("G828").Value = "Acelga"
$ws.Range("H828").Value = "Sin especificar"
$ws.Range("I828").Value = "Segunda"
$ws.Range("J828").Value = 60
$ws.Range("K828").Value = 12000
$ws.Range("L828").Value = 12000
$ws.Range("M828").Value = 12000
$ws.Range("N828").Value = "`$/docena de atados"
$ws.Range("O828").Value = "Región Metropolitana"
$ws.Range("P828").Value = 4000
$ws.Range("Q828").Value = 3
$ws.Range("R828").Value = "Hortaliza"
